$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19. This shifts the existing rows 19-55 down to 20-56,
# preserving all of their data/formatting (the new row inherits the date-format
# style from the row above, matching the target dimension A1:T56).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T carry over the same constant values used by
# every other row in this sheet (same market / product / variety / unit / origin).
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44622
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100104
$ws.Range("H19").Value = "Frutos de pepita"
$ws.Range("I19").Value = 100104003
$ws.Range("J19").Value = "Membrillo"
$ws.Range("K19").Value = "Champion"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 12
$ws.Range("N19").Value = 350000
$ws.Range("O19").Value = 360000
$ws.Range("P19").Value = 354167
$ws.Range("Q19").Value = "$/bins (450 kilos)"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 787
$ws.Range("T19").Value = 450
